$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "estado de cuenta" detail table (rows 16-20) previously listed 2 workers:
#   - CC 1047400500 INGRID TORRES TORRES  -> periods 2507, 2506, 2505 (rows 16-18)
#   - CC 1193123654 FABIAN ANDRES OñATE HERNANDEZ -> periods 2507, 2506 (rows 19-20)
#
# The update drops INGRID's records entirely and keeps only FABIAN's two periods,
# now listed with 2506 first and 2507 second (2507 keeps the bottom-border row
# style that used to sit on the 2506 row).
#
# Repurpose row 19 (currently FABIAN/2507) to hold FABIAN/2506's figures, and
# row 20 (currently FABIAN/2506, bottom-border style) to hold FABIAN/2507's
# figures, THEN delete the three INGRID rows (16-18) so rows 19/20 slide up
# into 16/17 - this keeps each row's original border styling attached to the
# correct (new) period.
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 51246

$ws.Range("E20").Value = "2507"
$ws.Range("F20").Value = 56940

$ws.Rows("16:18").Delete()

# Header summary figures: now 1 worker, 2 periods, and the total "Valor Mora"
# reflects just FABIAN's two remaining rows (51246 + 56940 = 108186).
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2
$ws.Range("E11").Value = 108186
